$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 408.66666
$ws.Range("I5").Value = 475.4
$ws.Range("J5").Value = 75
$ws.Range("K5").Value = 475.4
$ws.Range("L5").Value = 75
$ws.Range("M5").Value = -360.4
$ws.Range("N5").Value = -305
$ws.Range("H40").Value = 4258.0967
$ws.Range("J40").Value = 4511.222
$ws.Range("L40").Value = 4511.222
$ws.Range("N40").Value = -4861.222
$ws.Range("H69").Value = 8375
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 8773.4375
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 26320.3125
$ws.Range("M69").Value = -5126
$ws.Range("N69").Value = -28068.3125
$ws.Range("H72").Value = 8375
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 8773.4375
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 78960.9375
$ws.Range("M72").Value = -13632
$ws.Range("N72").Value = -87696.9375
$ws.Range("H107").Value = 2460.0356
$ws.Range("I107").Value = 2040.5
$ws.Range("K107").Value = 2040.5
$ws.Range("M107").Value = -120.5
$ws.Range("H141").Value = 1057.5
$ws.Range("I141").Value = 1057.5
$ws.Range("K141").Value = 3172.5
$ws.Range("M141").Value = 2007.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 1490.6
$ws.Range("I3").Value = 1490.6
$ws.Range("K3").Value = 1490.6
$ws.Range("M3").Value = -1375.6
$ws.Range("H11").Value = 10000000
$ws.Range("I11").Value = 10000000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 10000000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -9999856
$ws.Range("N11").ClearContents()
$ws.Range("H61").Value = 5565.1465
$ws.Range("I61").Value = 4835.1797
$ws.Range("K61").Value = 4835.1797
$ws.Range("M61").Value = -4623.1797
$ws.Range("H74").Value = 12347277
$ws.Range("I74").Value = 15153139
$ws.Range("J74").Value = 1482.8
$ws.Range("K74").Value = 15153139
$ws.Range("L74").Value = 1482.8
$ws.Range("M74").Value = -15152265
$ws.Range("N74").Value = -3230.8
$ws.Range("H77").Value = 12347277
$ws.Range("I77").Value = 15153139
$ws.Range("J77").Value = 1482.8
$ws.Range("K77").Value = 75765695
$ws.Range("L77").Value = 7414
$ws.Range("M77").Value = -75761327
$ws.Range("N77").Value = -16150
$ws.Range("H132").Value = 2383.9424
$ws.Range("I132").Value = 1876.2653
$ws.Range("K132").Value = 5628.7959
$ws.Range("M132").Value = -3098.7959
$ws.Range("H136").Value = 5565.1465
$ws.Range("I136").Value = 4835.1797
$ws.Range("K136").Value = 14505.5391
$ws.Range("M136").Value = -11955.5391

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 447.66666
$ws.Range("I22").Value = 337.2
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 337.2
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -164.2
$ws.Range("N22").Value = -1346
$ws.Range("H132").Value = 64998.6
$ws.Range("J132").Value = 64998.5
$ws.Range("L132").Value = 64998.5
$ws.Range("N132").Value = -75118.5
$ws.Range("H134").Value = 3435.9092
$ws.Range("I134").Value = 1987.25
$ws.Range("J134").Value = 7299
$ws.Range("K134").Value = 5961.75
$ws.Range("L134").Value = 21897
$ws.Range("M134").Value = -3426.75
$ws.Range("N134").Value = -26967

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 424.8889
$ws.Range("I7").Value = 520.5
$ws.Range("J7").Value = 348.4
$ws.Range("K7").Value = 520.5
$ws.Range("L7").Value = 348.4
$ws.Range("M7").Value = -407.5
$ws.Range("N7").Value = -574.4
$ws.Range("H22").Value = 1565.375
$ws.Range("I22").Value = 324.875
$ws.Range("J22").Value = 2805.875
$ws.Range("K22").Value = 324.875
$ws.Range("L22").Value = 2805.875
$ws.Range("M22").Value = 25.125
$ws.Range("N22").Value = -3505.875
$ws.Range("H31").Value = 29037.05
$ws.Range("I31").Value = 2798.037
$ws.Range("K31").Value = 2798.037
$ws.Range("M31").Value = -2503.037
$ws.Range("H34").Value = 29037.05
$ws.Range("I34").Value = 2798.037
$ws.Range("K34").Value = 2798.037
$ws.Range("M34").Value = -2596.037
$ws.Range("H58").Value = 3806.2693
$ws.Range("I58").Value = 2066.5557
$ws.Range("J58").Value = 7720.625
$ws.Range("K58").Value = 2066.5557
$ws.Range("L58").Value = 7720.625
$ws.Range("M58").Value = -1863.5557
$ws.Range("N58").Value = -8126.625
$ws.Range("H62").Value = 16253
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 16253
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 16253
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -17501
$ws.Range("H65").Value = 16253
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 16253
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 81265
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -87505
$ws.Range("H105").Value = 11249.4
$ws.Range("I105").Value = 7926.2856
$ws.Range("K105").Value = 7926.2856
$ws.Range("M105").Value = -6179.2856
$ws.Range("H134").Value = 3693.5833
$ws.Range("I134").Value = 1987
$ws.Range("J134").Value = 7106.75
$ws.Range("K134").Value = 5961
$ws.Range("L134").Value = 21320.25
$ws.Range("M134").Value = -3426
$ws.Range("N134").Value = -26390.25
$ws.Range("H136").Value = 3806.2693
$ws.Range("I136").Value = 2066.5557
$ws.Range("J136").Value = 7720.625
$ws.Range("K136").Value = 6199.6671
$ws.Range("L136").Value = 23161.875
$ws.Range("M136").Value = -3649.6671
$ws.Range("N136").Value = -28261.875
$ws.Range("H141").Value = 162213.67
$ws.Range("J141").Value = 162213.67
$ws.Range("L141").Value = 162213.67
$ws.Range("N141").Value = -172573.67
